# Add a new "functional_preprocessing.fmriInfo" parameter row to the
# functional pipeline parameter-properties worksheet.
#
# This corresponds to the commit:
#   "Add check on TR in default functional preprocessing script (#21)"
# which adds an fmriInfo option (used to adjust variables, e.g. TR, in the
# fMRI NIfTI header via mri_convert) right before the existing
# "functional_preprocessing.preprocessingScript" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 29
# ("functional_preprocessing.preprocessingScript"), which shifts that row
# and everything below it down by one.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new parameter's properties.
$ws.Cells.Item(29, 1).Value2 = "functional_preprocessing.fmriInfo"
$ws.Cells.Item(29, 4).Value2 = "functional_preprocessing"
$ws.Cells.Item(29, 5).Value2 = "char"
$ws.Cells.Item(29, 7).Value2 = "standard"
$ws.Cells.Item(29, 8).Value2 = 'Adjust variables in the fmriProcessedFile header (using mri_convert). Options are provided as structure (e.g. fmriInfo:{"tr": TR in msec, "te": TE in msec}). If emtpy, header is not changed.'

# Reflect the author's final cursor position in the saved workbook.
$ws.Range("D17").Select()
